$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16511686608523495"
$ws1.Range("B2").Value = "go_stims-16511686608105485.csv"
$ws1.Range("B3").Value = "GNG_stims-16511686608351176.csv"
$ws1.Range("B4").Value = "go_stims-165116866083712.csv"
$ws1.Range("B5").Value = "GNG_stims-16511686608503494.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-1651168663219381"
$ws2.Range("B2").Value = "TB-16511686631994596.csv"
$ws2.Range("B3").Value = "ZB-match_5-16511686613297598.csv"
$ws2.Range("B4").Value = "OB-1651168661463012.csv"
$ws2.Range("B5").Value = "OB-1651168661756671.csv"
$ws2.Range("B6").Value = "ZB-match_9-16511686610191736.csv"
$ws2.Range("B7").Value = "OB-16511686615539412.csv"
$ws2.Range("B8").Value = "TB-16511686622119231.csv"
$ws2.Range("B9").Value = "TB-16511686624758446.csv"
$ws2.Range("B10").Value = "ZB-match_2-16511686609024508.csv"

# --- Sheet 3: RS_TO (name only) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16511686632213836"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1651168663267382"
$ws4.Range("B2").Value = "MM_stims-16511686632352831.csv"
$ws4.Range("B3").Value = "ZM_stims-16511686632229013.csv"
$ws4.Range("B4").Value = "MM_stims-16511686632509606.csv"
$ws4.Range("B5").Value = "ZM_stims-16511686632352831.csv"
$ws4.Range("B6").Value = "MM_stims-16511686632663796.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686632519596.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16511686633441713"
$ws5.Range("B2").Value = "SAT_stims-16511686632976086.csv"
$ws5.Range("B3").Value = "vSAT_stims-16511686633289173.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651168663313744.csv"
$ws5.Range("B5").Value = "SAT_stims-16511686632713788.csv"
